$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newRows = @(
    @("2026/12/29","火",13,88),
    @("2026/12/29","火",16,99),
    @("2026/12/29","火",19,81),
    @("2026/12/29","火",23,85),
    @("2026/12/30","水",2,89),
    @("2026/12/30","水",5,95),
    @("2026/12/30","水",8,91),
    @("2026/12/30","水",13,92),
    @("2026/12/30","水",16,99),
    @("2026/12/30","水",22,108),
    @("2026/12/31","木",2,114),
    @("2026/12/31","木",6,120),
    @("2026/12/31","木",9,120),
    @("2026/12/31","木",12,201),
    @("2026/12/31","木",14,130),
    @("2026/12/31","木",22,120),
    @("2027/01/01","金",2,129),
    @("2027/01/01","金",5,119),
    @("2027/01/01","金",13,133),
    @("2027/01/01","金",16,109),
    @("2027/01/01","金",19,120),
    @("2027/01/02","土",1,105),
    @("2027/01/02","土",5,109),
    @("2027/01/02","土",8,110),
    @("2027/01/02","土",13,132),
    @("2027/01/02","土",16,145),
    @("2027/01/02","土",19,157),
    @("2027/01/02","土",22,165),
    @("2027/01/03","日",1,174),
    @("2027/01/03","日",4,192),
    @("2027/01/03","日",7,189),
    @("2027/01/03","日",13,201),
    @("2027/01/03","日",16,201),
    @("2027/01/03","日",19,201),
    @("2027/01/03","日",22,194),
    @("2027/01/04","月",2,164),
    @("2027/01/04","月",5,166),
    @("2027/01/04","月",7,168),
    @("2027/01/04","月",13,173),
    @("2027/01/04","月",22,127),
    @("2027/01/05","火",1,118),
    @("2027/01/05","火",7,127)
)

$startRow = 563
$endRow = 604

# Column A holds slash-separated dates as plain text (matches the rest of
# the sheet, which stores them as literal strings rather than Excel date
# serials) -- force text formatting before writing so COM doesn't coerce
# the "yyyy/mm/dd"-looking strings into date values.
$ws.Range("A$startRow`:A$endRow").NumberFormat = "@"

for ($i = 0; $i -lt $newRows.Count; $i++) {
    $r = $startRow + $i
    $row = $newRows[$i]
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Value = $row[3]
}
